# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 gets its table style switched from the
#    custom "Table_0" style to the built-in PowerPoint table style
#    {F8DB2291-9973-4511-B8E7-15183190F747}.
# 2) The presentation's applied theme colour scheme is changed from the
#    "Integral" / "Red Violet" palette back to the stock Office theme
#    palette ("Office").

function Hex2VbaRgb($hex) {
    # COM ColorFormat.RGB (and ThemeColorScheme item .RGB) use the
    # classic VBA RGB() packing: R + G*256 + B*65536.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1) Retarget the table's style on slide 5 -----------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{F8DB2291-9973-4511-B8E7-15183190F747}")
}

# --- 2) Restore the Office theme colour scheme -----------------------------
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.Slides.Item(1).Master.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = Hex2VbaRgb($officeColors[$i - 1])
}
